# Fixed error in excess mortality decomposition
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 5, 6, 7 (USA and extra rows no longer part of the table), and
# row 4 (old "Germany" row 21) - we will rewrite rows 2-4 entirely with the
# corrected data, then remove the now-unused rows 5-7.

# First, clear out rows 5 to 7 entirely.
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(5).Delete()

# Row 2: Germany, week 19 - corrected decomposition values
$ws.Range("A2").Value = "Germany"
$ws.Range("B2").Value = 19
$ws.Range("C2").Value = 0.03100350615901237
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").ClearContents()
$ws.Range("H2").ClearContents()

# Row 3: Spain, week 19
$ws.Range("A3").Value = "Spain"
$ws.Range("B3").Value = 19
$ws.Range("C3").Value = 0.1028601175810952
$ws.Range("D3").Value = -0.07185661142208284
$ws.Range("E3").Value = -0.03764720162676073
$ws.Range("F3").Value = -0.0342094097953221
$ws.Range("G3").Value = 0.5239211936341194
$ws.Range("H3").Value = 0.4760788063658806

# Row 4: Italy, week 19
$ws.Range("A4").Value = "Italy"
$ws.Range("B4").Value = 19
$ws.Range("C4").Value = 0.171962429438217
$ws.Range("D4").Value = -0.1409589232792046
$ws.Range("E4").Value = -0.06023535845635434
$ws.Range("F4").Value = -0.0807235648228503
$ws.Range("G4").Value = 0.427325614122655
$ws.Range("H4").Value = 0.5726743858773449
